$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1898.5
$ws.Range("J17").Value = 1898.5
$ws.Range("L17").Value = 5695.5
$ws.Range("N17").Value = -6031.5

$ws.Range("H70").Value = 4405.3
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 4405.3
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 13215.9
$ws.Range("N70").Value = -13755.9
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 4405.3
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 4405.3
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 13215.9
$ws.Range("N73").Value = -15087.9
$ws.Range("M73").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4423.706
$ws.Range("I61").Value = 2132.1428
$ws.Range("J61").Value = 6027.8
$ws.Range("K61").Value = 2132.1428
$ws.Range("L61").Value = 6027.8
$ws.Range("M61").Value = -1920.1428
$ws.Range("N61").Value = -6451.8

$ws.Range("H74").Value = 3469.0715
$ws.Range("I74").Value = 3505.6667
$ws.Range("J74").Value = 3249.5
$ws.Range("K74").Value = 3505.6667
$ws.Range("L74").Value = 3249.5
$ws.Range("M74").Value = -2631.6667
$ws.Range("N74").Value = -4997.5

$ws.Range("H77").Value = 3469.0715
$ws.Range("I77").Value = 3505.6667
$ws.Range("J77").Value = 3249.5
$ws.Range("K77").Value = 17528.3335
$ws.Range("L77").Value = 16247.5
$ws.Range("M77").Value = -13160.3335
$ws.Range("N77").Value = -24983.5

$ws.Range("H88").Value = 1250.0555
$ws.Range("J88").Value = 1088
$ws.Range("L88").Value = 1088
$ws.Range("N88").Value = -1900

$ws.Range("H91").Value = 1250.0555
$ws.Range("J91").Value = 1088
$ws.Range("L91").Value = 1088
$ws.Range("N91").Value = -3896

$ws.Range("H136").Value = 4423.706
$ws.Range("I136").Value = 2132.1428
$ws.Range("J136").Value = 6027.8
$ws.Range("K136").Value = 6396.428400000001
$ws.Range("L136").Value = 18083.4
$ws.Range("M136").Value = -3846.428400000001
$ws.Range("N136").Value = -23183.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3235.225
$ws.Range("I134").Value = 2394.3823
$ws.Range("K134").Value = 7183.146900000001
$ws.Range("M134").Value = -4648.146900000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1570.5714
$ws.Range("I4").Value = 1570.5714
$ws.Range("K4").Value = 1570.5714
$ws.Range("M4").Value = -1458.5714

$ws.Range("H7").Value = 201.94118
$ws.Range("I7").Value = 27.5
$ws.Range("J7").Value = 297.0909
$ws.Range("K7").Value = 27.5
$ws.Range("L7").Value = 297.0909
$ws.Range("M7").Value = 85.5
$ws.Range("N7").Value = -523.0908999999999

$ws.Range("H31").Value = 3278.1333
$ws.Range("I31").Value = 1384.2222
$ws.Range("K31").Value = 1384.2222
$ws.Range("M31").Value = -1089.2222

$ws.Range("H34").Value = 3278.1333
$ws.Range("I34").Value = 1384.2222
$ws.Range("K34").Value = 1384.2222
$ws.Range("M34").Value = -1182.2222

$ws.Range("H58").Value = 4714.3335
$ws.Range("I58").Value = 4912.25
$ws.Range("J58").Value = 4556
$ws.Range("K58").Value = 4912.25
$ws.Range("L58").Value = 4556
$ws.Range("M58").Value = -4709.25
$ws.Range("N58").Value = -4962

$ws.Range("H86").Value = 41608.625
$ws.Range("I86").Value = 54146.5
$ws.Range("K86").Value = 54146.5
$ws.Range("M86").Value = -53023.5

$ws.Range("H89").Value = 41608.625
$ws.Range("I89").Value = 54146.5
$ws.Range("K89").Value = 270732.5
$ws.Range("M89").Value = -265116.5

$ws.Range("H99").Value = 11817451
$ws.Range("J99").Value = 16674411
$ws.Range("L99").Value = 16674411
$ws.Range("N99").Value = -16677407

$ws.Range("H126").Value = 11817451
$ws.Range("J126").Value = 16674411
$ws.Range("L126").Value = 50023233
$ws.Range("N126").Value = -50028173

$ws.Range("H132").Value = 2988
$ws.Range("I132").Value = 2986.7896
$ws.Range("K132").Value = 8960.3688
$ws.Range("M132").Value = -6430.3688

$ws.Range("H134").Value = 4613.3228
$ws.Range("I134").Value = 3426.5264
$ws.Range("K134").Value = 10279.5792
$ws.Range("M134").Value = -7744.5792

$ws.Range("H136").Value = 4714.3335
$ws.Range("I136").Value = 4912.25
$ws.Range("J136").Value = 4556
$ws.Range("K136").Value = 14736.75
$ws.Range("L136").Value = 13668
$ws.Range("M136").Value = -12186.75
$ws.Range("N136").Value = -18768

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3496.5
$ws.Range("I3").Value = 3496.5
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 10489.5
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -10377.5
$ws.Range("N3").ClearContents()

$ws.Range("H107").Value = 566.6667
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 126801.11
$ws.Range("I80").Value = 187201.83
$ws.Range("K80").Value = 187201.83
$ws.Range("M80").Value = -186203.83

$ws.Range("H83").Value = 126801.11
$ws.Range("I83").Value = 187201.83
$ws.Range("K83").Value = 936009.1499999999
$ws.Range("M83").Value = -931017.1499999999

$ws.Range("H102").Value = 1450
$ws.Range("I102").Value = 1170
$ws.Range("J102").Value = 1800
$ws.Range("K102").Value = 1170
$ws.Range("L102").Value = 1800
$ws.Range("M102").Value = 452
$ws.Range("N102").Value = -5044

$ws.Range("H132").Value = 4237.857
$ws.Range("I132").Value = 2608.375
$ws.Range("J132").Value = 6410.5
$ws.Range("K132").Value = 7825.125
$ws.Range("L132").Value = 19231.5
$ws.Range("M132").Value = -5295.125
$ws.Range("N132").Value = -24291.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3626
$ws.Range("I7").Value = 2405.3333
$ws.Range("J7").Value = 5290.5454
$ws.Range("K7").Value = 2405.3333
$ws.Range("L7").Value = 5290.5454
$ws.Range("M7").Value = -2293.3333
$ws.Range("N7").Value = -5514.5454

$ws.Range("H40").Value = 16961.2
$ws.Range("J40").Value = 9142.857
$ws.Range("L40").Value = 9142.857
$ws.Range("N40").Value = -9414.857

$ws.Range("H46").Value = 1835.8462
$ws.Range("I46").Value = 994.5
$ws.Range("K46").Value = 994.5
$ws.Range("M46").Value = -806.5

$ws.Range("H101").Value = 31103.285
$ws.Range("J101").Value = 31103.285
$ws.Range("L101").Value = 31103.285
$ws.Range("N101").Value = -37593.285

$ws.Range("H126").Value = 3626
$ws.Range("I126").Value = 2405.3333
$ws.Range("J126").Value = 5290.5454
$ws.Range("K126").Value = 7215.999899999999
$ws.Range("L126").Value = 15871.6362
$ws.Range("M126").Value = -4745.999899999999
$ws.Range("N126").Value = -20811.6362

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 7900
$ws.Range("J5").Value = 7900
$ws.Range("L5").Value = 7900
$ws.Range("N5").Value = -8124

$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()

$ws.Range("H55").Value = 4950
$ws.Range("I55").Value = 4900
$ws.Range("K55").Value = 4900
$ws.Range("M55").Value = -4623

$ws.Range("H59").Value = 16000.5
$ws.Range("J59").Value = 16000.5
$ws.Range("L59").Value = 16000.5
$ws.Range("N59").Value = -17476.5

$ws.Range("H61").Value = 19271.182
$ws.Range("I61").Value = 16748.8
$ws.Range("K61").Value = 16748.8
$ws.Range("M61").Value = -16456.8

$ws.Range("H81").Value = 8366.200000000001
$ws.Range("J81").Value = 1805.3334
$ws.Range("L81").Value = 3610.6668
$ws.Range("N81").Value = -5732.6668

$ws.Range("H84").Value = 8366.200000000001
$ws.Range("J84").Value = 1805.3334
$ws.Range("L84").Value = 18053.334
$ws.Range("N84").Value = -28661.334

$ws.Range("H136").Value = 3022.2812
$ws.Range("I136").Value = 1807.5454
$ws.Range("K136").Value = 5422.6362
$ws.Range("M136").Value = -2872.6362
